$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 97.57717366666668
$ws.Range("H2").Value = 292.731521
$ws.Range("I2").Value = 0.3532166605548384
$ws.Range("J2").Value = 0.3532166605548384
$ws.Range("M2").Value = 1.231278
$ws.Range("N2").Value = 3.693834
$ws.Range("O2").Value = 0.3283785416403858
$ws.Range("P2").Value = 0.3283785416403859
$ws.Range("Q2").Value = 120.144627237946
$ws.Range("R2").Value = 1081.301645141514
$ws.Range("S2").Value = 0.115988771876085
$ws.Range("T2").Value = 0.1159887718760851
$ws.Range("G3").Value = 97.57717366666668
$ws.Range("H3").Value = 292.731521
$ws.Range("I3").Value = 0.3532166605548384
$ws.Range("J3").Value = 0.3532166605548384
$ws.Range("O3").Value = 0.2084514246837437
$ws.Range("P3").Value = 0.2084514246837437
$ws.Range("Q3").Value = 76.26661166938769
$ws.Range("R3").Value = 686.3995050244891
$ws.Range("S3").Value = 0.07362851611469035
$ws.Range("T3").Value = 0.07362851611469036
$ws.Range("G4").Value = 97.57717366666668
$ws.Range("H4").Value = 292.731521
$ws.Range("I4").Value = 0.3532166605548384
$ws.Range("J4").Value = 0.3532166605548384
$ws.Range("O4").Value = 0.4631700336758705
$ws.Range("P4").Value = 0.4631700336758705
$ws.Range("Q4").Value = 169.461106580816
$ws.Range("R4").Value = 1525.149959227344
$ws.Range("S4").Value = 0.163599372564063
$ws.Range("T4").Value = 0.163599372564063
$ws.Range("I5").Value = 0.5533024543641269
$ws.Range("J5").Value = 0.5533024543641269
$ws.Range("M5").Value = 1.231278
$ws.Range("N5").Value = 3.693834
$ws.Range("O5").Value = 0.3283785416403858
$ws.Range("P5").Value = 0.3283785416403859
$ws.Range("Q5").Value = 188.202665822718
$ws.Range("R5").Value = 1693.823992404462
$ws.Range("S5").Value = 0.1816926530501381
$ws.Range("T5").Value = 0.1816926530501382
$ws.Range("I6").Value = 0.5533024543641269
$ws.Range("J6").Value = 0.5533024543641269
$ws.Range("O6").Value = 0.2084514246837437
$ws.Range("P6").Value = 0.2084514246837437
$ws.Range("S6").Value = 0.1153366848932143
$ws.Range("T6").Value = 0.1153366848932143
$ws.Range("I7").Value = 0.5533024543641269
$ws.Range("J7").Value = 0.5533024543641269
$ws.Range("O7").Value = 0.4631700336758705
$ws.Range("P7").Value = 0.4631700336758705
$ws.Range("S7").Value = 0.2562731164207744
$ws.Range("T7").Value = 0.2562731164207745
$ws.Range("I8").Value = 0.09348088508103472
$ws.Range("J8").Value = 0.09348088508103473
$ws.Range("M8").Value = 1.231278
$ws.Range("N8").Value = 3.693834
$ws.Range("O8").Value = 0.3283785416403858
$ws.Range("P8").Value = 0.3283785416403859
$ws.Range("Q8").Value = 31.796988495084
$ws.Range("R8").Value = 286.172896455756
$ws.Range("S8").Value = 0.03069711671416268
$ws.Range("T8").Value = 0.03069711671416269
$ws.Range("I9").Value = 0.09348088508103472
$ws.Range("J9").Value = 0.09348088508103473
$ws.Range("O9").Value = 0.2084514246837437
$ws.Range("P9").Value = 0.2084514246837437
$ws.Range("S9").Value = 0.01948622367583901
$ws.Range("T9").Value = 0.01948622367583901
$ws.Range("I10").Value = 0.09348088508103472
$ws.Range("J10").Value = 0.09348088508103473
$ws.Range("O10").Value = 0.4631700336758705
$ws.Range("P10").Value = 0.4631700336758705
$ws.Range("S10").Value = 0.04329754469103303
$ws.Range("T10").Value = 0.04329754469103304
